$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1060.05
$ws.Range("I19").Value = 999.9231
$ws.Range("J19").Value = 1171.7142
$ws.Range("K19").Value = 999.9231
$ws.Range("L19").Value = 1171.7142
$ws.Range("M19").Value = -824.9231
$ws.Range("N19").Value = -1521.7142
$ws.Range("H34").Value = 12117.625
$ws.Range("I34").Value = 12117.625
$ws.Range("K34").Value = 12117.625
$ws.Range("M34").Value = -11914.625
$ws.Range("H36").Value = 12117.625
$ws.Range("I36").Value = 12117.625
$ws.Range("K36").Value = 12117.625
$ws.Range("M36").Value = -11402.625
$ws.Range("H43").Value = 5269.5
$ws.Range("J43").Value = 3372.875
$ws.Range("L43").Value = 3372.875
$ws.Range("N43").Value = -3510.875
$ws.Range("H53").Value = 874.0769
$ws.Range("I53").Value = 370.22223
$ws.Range("K53").Value = 370.22223
$ws.Range("M53").Value = 266.77777
$ws.Range("H62").Value = 16002
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 16002
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 16002
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -17250
$ws.Range("H65").Value = 16002
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 16002
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 80010
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -86250
$ws.Range("H88").Value = 766.6923
$ws.Range("J88").Value = 820.36365
$ws.Range("L88").Value = 820.36365
$ws.Range("N88").Value = -1632.36365
$ws.Range("H91").Value = 766.6923
$ws.Range("J91").Value = 820.36365
$ws.Range("L91").Value = 820.36365
$ws.Range("N91").Value = -3628.36365
$ws.Range("H116").Value = 7396.6665
$ws.Range("J116").Value = 6625.6665
$ws.Range("L116").Value = 6625.6665
$ws.Range("N116").Value = -13509.6665
$ws.Range("H127").Value = 1498.3914
$ws.Range("I127").Value = 971.7368
$ws.Range("K127").Value = 2915.2104
$ws.Range("M127").Value = 2044.7896
$ws.Range("H129").Value = 7080.727
$ws.Range("I129").Value = 7288.9
$ws.Range("K129").Value = 21866.7
$ws.Range("M129").Value = -16866.7
$ws.Range("H135").Value = 1070.8572
$ws.Range("I135").Value = 1134
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 10206
$ws.Range("L135").Value = 2250
$ws.Range("M135").Value = -7671
$ws.Range("N135").Value = -7320
$ws.Range("H137").Value = 2022.0233
$ws.Range("I137").Value = 2056.4
$ws.Range("J137").Value = 1563.6666
$ws.Range("K137").Value = 6169.200000000001
$ws.Range("L137").Value = 4690.9998
$ws.Range("M137").Value = -3619.200000000001
$ws.Range("N137").Value = -9790.9998
$ws.Range("H138").Value = 5464.0723
$ws.Range("I138").Value = 4236.255
$ws.Range("J138").Value = 7420.9062
$ws.Range("K138").Value = 12708.765
$ws.Range("L138").Value = 22262.7186
$ws.Range("M138").Value = -7568.764999999999
$ws.Range("N138").Value = -32542.7186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3321.9092
$ws.Range("I45").Value = 1658.75
$ws.Range("K45").Value = 1658.75
$ws.Range("M45").Value = -1281.75
$ws.Range("H74").Value = 45458772
$ws.Range("I74").Value = 45458772
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 45458772
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -45457898
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 45458772
$ws.Range("I77").Value = 45458772
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 227293860
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -227289492
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 2242.0454
$ws.Range("I110").Value = 2353.3157
$ws.Range("K110").Value = 2353.3157
$ws.Range("M110").Value = -308.3157000000001
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800
$ws.Range("H132").Value = 3431.465
$ws.Range("I132").Value = 1986.6
$ws.Range("J132").Value = 9752.75
$ws.Range("K132").Value = 5959.799999999999
$ws.Range("L132").Value = 29258.25
$ws.Range("M132").Value = -3429.799999999999
$ws.Range("N132").Value = -34318.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1681.5428
$ws.Range("I20").Value = 1537.5238
$ws.Range("J20").Value = 1897.5714
$ws.Range("K20").Value = 1537.5238
$ws.Range("L20").Value = 1897.5714
$ws.Range("M20").Value = -1290.5238
$ws.Range("N20").Value = -2391.5714
$ws.Range("H21").Value = 36660
$ws.Range("J21").Value = 36660
$ws.Range("L21").Value = 36660
$ws.Range("N21").Value = -37132
$ws.Range("H22").Value = 567.35
$ws.Range("I22").Value = 594.05554
$ws.Range("K22").Value = 594.05554
$ws.Range("M22").Value = -421.05554
$ws.Range("H29").Value = 1046
$ws.Range("I29").Value = 1046
$ws.Range("K29").Value = 1046
$ws.Range("M29").Value = -757
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H94").Value = 598.9167
$ws.Range("I94").Value = 598.9167
$ws.Range("K94").Value = 598.9167
$ws.Range("M94").Value = -147.9167
$ws.Range("H96").Value = 15644.25
$ws.Range("I96").Value = 11150.571
$ws.Range("K96").Value = 11150.571
$ws.Range("M96").Value = -8404.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 7250
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 7250
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 7250
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -7836
$ws.Range("H31").Value = 7748.2856
$ws.Range("I31").Value = 2717.7827
$ws.Range("K31").Value = 2717.7827
$ws.Range("M31").Value = -2422.7827
$ws.Range("H34").Value = 7748.2856
$ws.Range("I34").Value = 2717.7827
$ws.Range("K34").Value = 2717.7827
$ws.Range("M34").Value = -2515.7827
$ws.Range("H86").Value = 4284.7
$ws.Range("I86").Value = 3975
$ws.Range("J86").Value = 4362.125
$ws.Range("K86").Value = 3975
$ws.Range("L86").Value = 4362.125
$ws.Range("M86").Value = -2852
$ws.Range("N86").Value = -6608.125
$ws.Range("H89").Value = 4284.7
$ws.Range("I89").Value = 3975
$ws.Range("J89").Value = 4362.125
$ws.Range("K89").Value = 19875
$ws.Range("L89").Value = 21810.625
$ws.Range("M89").Value = -14259
$ws.Range("N89").Value = -33042.625
$ws.Range("H105").Value = 2765.8333
$ws.Range("I105").Value = 3019
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 3019
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -1272
$ws.Range("N105").Value = -4994
$ws.Range("H122").Value = 2818.5557
$ws.Range("I122").Value = 2094.6667
$ws.Range("K122").Value = 6284.000100000001
$ws.Range("M122").Value = -3834.000100000001
$ws.Range("H124").Value = 7044975.5
$ws.Range("J124").Value = 7044975.5
$ws.Range("L124").Value = 7044975.5
$ws.Range("N124").Value = -7049885.5
$ws.Range("H132").Value = 2666.484
$ws.Range("I132").Value = 1609.0344
$ws.Range("K132").Value = 4827.1032
$ws.Range("M132").Value = -2297.1032

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1496.3077
$ws.Range("I5").Value = 1780
$ws.Range("J5").Value = 1253.1428
$ws.Range("K5").Value = 5340
$ws.Range("L5").Value = 3759.4284
$ws.Range("M5").Value = -5228
$ws.Range("N5").Value = -3983.4284
$ws.Range("H25").Value = 100
$ws.Range("I25").Value = 100
$ws.Range("K25").Value = 300
$ws.Range("M25").Value = -131
$ws.Range("H30").Value = 100
$ws.Range("I30").Value = 100
$ws.Range("K30").Value = 300
$ws.Range("M30").Value = -198
$ws.Range("H36").Value = 450
$ws.Range("I36").Value = 450
$ws.Range("K36").Value = 1350
$ws.Range("M36").Value = -1181
$ws.Range("H75").Value = 4378.8887
$ws.Range("J75").Value = 5488.4287
$ws.Range("L75").Value = 16465.2861
$ws.Range("N75").Value = -18461.2861
$ws.Range("H78").Value = 4378.8887
$ws.Range("J78").Value = 5488.4287
$ws.Range("L78").Value = 49395.85830000001
$ws.Range("N78").Value = -59379.85830000001
$ws.Range("H135").Value = 1496.3077
$ws.Range("I135").Value = 1780
$ws.Range("J135").Value = 1253.1428
$ws.Range("K135").Value = 16020
$ws.Range("L135").Value = 11278.2852
$ws.Range("M135").Value = -13485
$ws.Range("N135").Value = -16348.2852
$ws.Range("H137").Value = 3851.652
$ws.Range("I137").Value = 2906.8333
$ws.Range("J137").Value = 4185.1177
$ws.Range("K137").Value = 8720.499899999999
$ws.Range("L137").Value = 12555.3531
$ws.Range("M137").Value = -3620.499899999999
$ws.Range("N137").Value = -22755.3531
$ws.Range("H140").Value = 2024.8379
$ws.Range("I140").Value = 1962.5
$ws.Range("J140").Value = 2047.9259
$ws.Range("K140").Value = 5887.5
$ws.Range("L140").Value = 6143.7777
$ws.Range("M140").Value = -707.5
$ws.Range("N140").Value = -16503.7777

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6832.2607
$ws.Range("I70").Value = 5743.75
$ws.Range("K70").Value = 5743.75
$ws.Range("M70").Value = -5473.75
$ws.Range("H73").Value = 6832.2607
$ws.Range("I73").Value = 5743.75
$ws.Range("K73").Value = 5743.75
$ws.Range("M73").Value = -4807.75
$ws.Range("H80").Value = 5498.8213
$ws.Range("I80").Value = 3866.3684
$ws.Range("K80").Value = 3866.3684
$ws.Range("M80").Value = -2868.3684
$ws.Range("H82").Value = 59999
$ws.Range("J82").Value = 59999
$ws.Range("L82").Value = 59999
$ws.Range("N82").Value = -60765
$ws.Range("H83").Value = 5498.8213
$ws.Range("I83").Value = 3866.3684
$ws.Range("K83").Value = 19331.842
$ws.Range("M83").Value = -14339.842
$ws.Range("H85").Value = 59999
$ws.Range("J85").Value = 59999
$ws.Range("L85").Value = 59999
$ws.Range("N85").Value = -62651
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 120000
$ws.Range("N136").Value = -125100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3753.8333
$ws.Range("I7").Value = 3790.739
$ws.Range("K7").Value = 3790.739
$ws.Range("M7").Value = -3678.739
$ws.Range("H16").Value = 2659.7307
$ws.Range("I16").Value = 497.52942
$ws.Range("J16").Value = 6743.8887
$ws.Range("K16").Value = 497.52942
$ws.Range("L16").Value = 6743.8887
$ws.Range("M16").Value = -327.52942
$ws.Range("N16").Value = -7083.8887
$ws.Range("H22").Value = 1708.6
$ws.Range("I22").Value = 1189.4286
$ws.Range("K22").Value = 1189.4286
$ws.Range("M22").Value = -894.4286
$ws.Range("H27").Value = 1708.6
$ws.Range("I27").Value = 1189.4286
$ws.Range("K27").Value = 1189.4286
$ws.Range("M27").Value = -1082.4286
$ws.Range("H42").Value = 62000
$ws.Range("J42").Value = 62000
$ws.Range("L42").Value = 62000
$ws.Range("N42").Value = -63126
$ws.Range("H46").Value = 2523.4
$ws.Range("I46").Value = 845.8
$ws.Range("J46").Value = 4201
$ws.Range("K46").Value = 845.8
$ws.Range("L46").Value = 4201
$ws.Range("M46").Value = -657.8
$ws.Range("N46").Value = -4577
$ws.Range("H49").Value = 62000
$ws.Range("J49").Value = 62000
$ws.Range("L49").Value = 62000
$ws.Range("N49").Value = -62294
$ws.Range("H55").Value = 814.6667
$ws.Range("I55").Value = 147.7
$ws.Range("K55").Value = 147.7
$ws.Range("M55").Value = 25.30000000000001
$ws.Range("H69").Value = 67388
$ws.Range("J69").Value = 67388
$ws.Range("L69").Value = 67388
$ws.Range("N69").Value = -69010
$ws.Range("H72").Value = 67388
$ws.Range("J72").Value = 67388
$ws.Range("L72").Value = 202164
$ws.Range("N72").Value = -210276
$ws.Range("H82").Value = 1197.1538
$ws.Range("I82").Value = 706
$ws.Range("J82").Value = 1770.1666
$ws.Range("K82").Value = 706
$ws.Range("L82").Value = 1770.1666
$ws.Range("M82").Value = -345
$ws.Range("N82").Value = -2492.1666
$ws.Range("H85").Value = 1197.1538
$ws.Range("I85").Value = 706
$ws.Range("J85").Value = 1770.1666
$ws.Range("K85").Value = 706
$ws.Range("L85").Value = 1770.1666
$ws.Range("M85").Value = 542
$ws.Range("N85").Value = -4266.1666
$ws.Range("H93").Value = 978.26666
$ws.Range("I93").Value = 882.6539
$ws.Range("K93").Value = 882.6539
$ws.Range("M93").Value = 365.3461
$ws.Range("H100").Value = 4277.857
$ws.Range("I100").Value = 7749
$ws.Range("J100").Value = 2889.4
$ws.Range("K100").Value = 7749
$ws.Range("L100").Value = 2889.4
$ws.Range("M100").Value = -7208
$ws.Range("N100").Value = -3971.4
$ws.Range("H122").Value = 4638.926
$ws.Range("I122").Value = 3492.2104
$ws.Range("J122").Value = 7362.375
$ws.Range("K122").Value = 10476.6312
$ws.Range("L122").Value = 22087.125
$ws.Range("M122").Value = -8026.6312
$ws.Range("N122").Value = -26987.125
$ws.Range("H126").Value = 3753.8333
$ws.Range("I126").Value = 3790.739
$ws.Range("K126").Value = 11372.217
$ws.Range("M126").Value = -8902.217000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 875.6667
$ws.Range("I113").Value = 782.2308
$ws.Range("K113").Value = 2346.6924
$ws.Range("M113").Value = -176.6923999999999
